$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 22:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1682428
$ws.Range("C4").Value = 15600
$ws.Range("D4").Value = 451392
$ws.Range("E4").Value = 1131811
$ws.Range("G4").Value = 542
$ws.Range("H4").Value = 99225

# Row 5 - Brasil
$ws.Range("B5").Value = 357839
$ws.Range("C5").Value = 10441
$ws.Range("E5").Value = 192752
$ws.Range("G5").Value = 487
$ws.Range("H5").Value = 22500

# Row 11 - Alemania
$ws.Range("B11").Value = 180321
$ws.Range("C11").Value = 335
$ws.Range("E11").Value = 11650

# Row 15 - Peru
$ws.Range("B15").Value = 119959
$ws.Range("C15").Value = 4205
$ws.Range("E15").Value = 68671

# Row 84 - Costa de Marfil
$ws.Range("B84").Value = 2376
$ws.Range("C84").Value = 10
$ws.Range("D84").Value = 1219
$ws.Range("E84").Value = 1127

# Row 152 - Mauritania
$ws.Range("B152").Value = 237
$ws.Range("C152").Value = 10
$ws.Range("E152").Value = 216
